$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.403.92"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "2.076.10"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.30"
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.624"
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.12"
$ws.Range("E8").Value = "  -2.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.383"
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0764"
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").Value = "2.373.75"
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.66"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.78"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.779"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.15"
$ws.Range("E16").Value = "  -2.16%  "
$ws.Range("D17").Value = "2.068.10"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").Value = "37.323.78"
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.31"
$ws.Range("E19").Value = "  +4.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.53"
$ws.Range("E20").Value = "  +1.28%  "
$ws.Range("D21").Value = "0.0₃0813"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.69"
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.40"
$ws.Range("E25").Value = "  -2.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.67"
$ws.Range("E26").Value = "  +2.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.79"
$ws.Range("E27").Value = "  -1.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.44"
$ws.Range("E28").Value = "  +3.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.11"
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("E30").Value = "  -4.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.118"
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.46"
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0618"
$ws.Range("E33").Value = "  -1.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.57"
$ws.Range("E34").Value = "  +2.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.50"
$ws.Range("E35").Value = "  -4.76%  "
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.76"
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.25"
$ws.Range("E38").Value = "  -3.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.67"
$ws.Range("E39").Value = "  -4.95%  "
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("B41").Value = "FTXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.43"
$ws.Range("E41").Value = "  +3.18%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.473.64"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "96.44"
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0941"
$ws.Range("E44").Value = "  -3.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0213"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.17"
$ws.Range("E46").Value = "  +2.87%  "
$ws.Range("E47").Value = "  -0.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.13"
$ws.Range("E48").Value = "  -8.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.16"
$ws.Range("E49").Value = "  -2.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.97"
$ws.Range("E50").Value = "  +1.04%  "
$ws.Range("D51").Value = "2.261.15"
$ws.Range("E51").Value = "  -0.30%  "
